$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 61, shifting elz99 (currently row 61) down to row 62
$ws.Rows.Item(61).Insert()

# Fill in the new row 61 with elz93 / 99993
$ws.Range("A61").Value = "elz93"
$ws.Range("B61").Value = 99993

# Copy the style of B62 (elz99's row, style index 23) onto the new B61 cell
$ws.Range("B62").Copy()
$ws.Range("B61").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
